{"js": "// Update the worksheet date and the 25 division-problem answers in the\n// 5x5 grid (stored as rows 0,4,8,12,16 of the single table; the other\n// rows are blank spacer rows).\n\n// 1) Update the date/weekday line (first paragraph of the body, outside\n//    the table).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\nif (dateParagraph.text === \"2023-10-22 Sunday\") {\n  dateParagraph.getRange().insertText(\"2023-10-23 Monday\", \"Replace\");\n}\n\n// 2) Update the 25 division problems inside the table. Cells are\n//    addressed by (row, col) on the full table grid (0-based), which\n//    includes the blank spacer rows between each data row.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst dataRows = [0, 4, 8, 12, 16];\nconst oldGrid = [\n  [\"86\u00f79=9, 5\", \"85\u00f77=12, 1\", \"12\u00f75=2, 2\", \"81\u00f72=40, 1\", \"21\u00f75=4, 1\"],\n  [\"87\u00f73=29, 0\", \"93\u00f75=18, 3\", \"26\u00f77=3, 5\", \"36\u00f73=12, 0\", \"72\u00f78=9, 0\"],\n  [\"77\u00f74=19, 1\", \"16\u00f78=2, 0\", \"19\u00f78=2, 3\", \"18\u00f74=4, 2\", \"69\u00f75=13, 4\"],\n  [\"67\u00f72=33, 1\", \"54\u00f75=10, 4\", \"75\u00f72=37, 1\", \"94\u00f72=47, 0\", \"60\u00f75=12, 0\"],\n  [\"41\u00f79=4, 5\", \"30\u00f72=15, 0\", \"67\u00f72=33, 1\", \"50\u00f72=25, 0\", \"36\u00f75=7, 1\"],\n];\nconst newGrid = [\n  [\"86\u00f74=21, 2\", \"55\u00f78=6, 7\", \"87\u00f79=9, 6\", \"55\u00f76=9, 1\", \"53\u00f77=7, 4\"],\n  [\"47\u00f79=5, 2\", \"95\u00f77=13, 4\", \"22\u00f74=5, 2\", \"73\u00f72=36, 1\", \"65\u00f75=13, 0\"],\n  [\"52\u00f77=7, 3\", \"17\u00f77=2, 3\", \"85\u00f78=10, 5\", \"80\u00f72=40, 0\", \"40\u00f74=10, 0\"],\n  [\"29\u00f72=14, 1\", \"15\u00f74=3, 3\", \"98\u00f73=32, 2\", \"24\u00f77=3, 3\", \"12\u00f78=1, 4\"],\n  [\"11\u00f73=3, 2\", \"77\u00f79=8, 5\", \"31\u00f72=15, 1\", \"36\u00f77=5, 1\", \"83\u00f77=11, 6\"],\n];\n\nfor (let r = 0; r < dataRows.length; r++) {\n  const tableRow = dataRows[r];\n  for (let c = 0; c < oldGrid[r].length; c++) {\n    const cell = table.getCell(tableRow, c);\n    cell.load(\"value\");\n    await context.sync();\n    if (cell.value === oldGrid[r][c]) {\n      cell.value = newGrid[r][c];\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division-problem answers in the\n# 5x5 grid (stored as rows 1,5,9,13,17 of the single table, 1-based;\n# the other rows are blank spacer rows).\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/weekday line (first paragraph of the body, outside\n#    the table). Paragraph.Range.Text carries the trailing paragraph\n#    mark (CR, char 13), so trim it before comparing.\n$p1 = $d.Paragraphs.Item(1)\n$p1Text = $p1.Range.Text.TrimEnd([char]13)\nif ($p1Text -eq \"2023-10-22 Sunday\") {\n    $p1.Range.Text = \"2023-10-23 Monday\"\n}\n\n# 2) Update the 25 division problems inside the table. Cells are\n#    addressed by (row, col), 1-based, on the full table grid, which\n#    includes the blank spacer rows between each data row.\n$table = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\n$oldGrid = @(\n    @(\"86\u00f79=9, 5\", \"85\u00f77=12, 1\", \"12\u00f75=2, 2\", \"81\u00f72=40, 1\", \"21\u00f75=4, 1\"),\n    @(\"87\u00f73=29, 0\", \"93\u00f75=18, 3\", \"26\u00f77=3, 5\", \"36\u00f73=12, 0\", \"72\u00f78=9, 0\"),\n    @(\"77\u00f74=19, 1\", \"16\u00f78=2, 0\", \"19\u00f78=2, 3\", \"18\u00f74=4, 2\", \"69\u00f75=13, 4\"),\n    @(\"67\u00f72=33, 1\", \"54\u00f75=10, 4\", \"75\u00f72=37, 1\", \"94\u00f72=47, 0\", \"60\u00f75=12, 0\"),\n    @(\"41\u00f79=4, 5\", \"30\u00f72=15, 0\", \"67\u00f72=33, 1\", \"50\u00f72=25, 0\", \"36\u00f75=7, 1\")\n)\n\n$newGrid = @(\n    @(\"86\u00f74=21, 2\", \"55\u00f78=6, 7\", \"87\u00f79=9, 6\", \"55\u00f76=9, 1\", \"53\u00f77=7, 4\"),\n    @(\"47\u00f79=5, 2\", \"95\u00f77=13, 4\", \"22\u00f74=5, 2\", \"73\u00f72=36, 1\", \"65\u00f75=13, 0\"),\n    @(\"52\u00f77=7, 3\", \"17\u00f77=2, 3\", \"85\u00f78=10, 5\", \"80\u00f72=40, 0\", \"40\u00f74=10, 0\"),\n    @(\"29\u00f72=14, 1\", \"15\u00f74=3, 3\", \"98\u00f73=32, 2\", \"24\u00f77=3, 3\", \"12\u00f78=1, 4\"),\n    @(\"11\u00f73=3, 2\", \"77\u00f79=8, 5\", \"31\u00f72=15, 1\", \"36\u00f77=5, 1\", \"83\u00f77=11, 6\")\n)\n\nfor ($r = 0; $r -lt $dataRows.Length; $r++) {\n    $tableRow = $dataRows[$r]\n    for ($c = 0; $c -lt 5; $c++) {\n        $cell = $table.Cell($tableRow, $c + 1)\n        $old = $oldGrid[$r][$c]\n        $new = $newGrid[$r][$c]\n        # Cell.Range.Text carries a trailing end-of-cell mark (CR + BEL,\n        # chars 13/7) that isn't part of the visible text, so trim it\n        # before comparing.\n        $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($current -eq $old) {\n            $cell.Range.Text = $new\n        }\n    }\n}\n"}
